$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.921622333333333
$ws.Range("N2").Value = 5.764867
$ws.Range("O2").Value = 0.1392241219313625
$ws.Range("P2").Value = 0.1392241219313625
$ws.Range("Q2").Value = 0.138696935153
$ws.Range("R2").Value = 1.248272416377
$ws.Range("S2").Value = 0.1392241219313625
$ws.Range("T2").Value = 0.1392241219313625

# Row 3
$ws.Range("O3").Value = 0.7511588049189343
$ws.Range("P3").Value = 0.7511588049189343
$ws.Range("S3").Value = 0.7511588049189343
$ws.Range("T3").Value = 0.7511588049189343

# Row 4
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.1441973333333333
$ws.Range("N4").Value = 0.432592
$ws.Range("O4").Value = 0.01044729069283506
$ws.Range("P4").Value = 0.01044729069283506
$ws.Range("Q4").Value = 0.010407730928
$ws.Range("R4").Value = 0.093669578352
$ws.Range("S4").Value = 0.01044729069283506
$ws.Range("T4").Value = 0.01044729069283506

# Row 5
$ws.Range("M5").Value = 1.182384
$ws.Range("N5").Value = 3.547152
$ws.Range("O5").Value = 0.08566531067535062
$ws.Range("P5").Value = 0.08566531067535062
$ws.Range("Q5").Value = 0.08534092996800001
$ws.Range("R5").Value = 0.768068369712
$ws.Range("S5").Value = 0.08566531067535062
$ws.Range("T5").Value = 0.08566531067535062

# Row 6
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1863936666666667
$ws.Range("N6").Value = 0.559181
$ws.Range("O6").Value = 0.01350447178151746
$ws.Range("P6").Value = 0.01350447178151746
$ws.Range("Q6").Value = 0.013453335679
$ws.Range("R6").Value = 0.121080021111
$ws.Range("S6").Value = 0.01350447178151746
$ws.Range("T6").Value = 0.01350447178151746
